$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.06257505313200211
$ws.Range("H2").Value = -2.685984518905973
$ws.Range("I2").Value = -6.302553151364791
$ws.Range("G3").Value = 0.06811769916351476
$ws.Range("H3").Value = 21.19068956614814
$ws.Range("G4").Value = -0.02276256339526492
$ws.Range("H4").Value = -12.06677159169732
$ws.Range("G5").Value = -0.0224435486643241
$ws.Range("H5").Value = -96.33364421389865
$ws.Range("G6").Value = -0.0101449860531542
$ws.Range("H6").Value = 9.367634910516596
$ws.Range("G7").Value = -0.02020051505303744
$ws.Range("H7").Value = -256.6641350000703
$ws.Range("G8").Value = -0.0008845633449788281
$ws.Range("H8").Value = 84.4517841127785
$ws.Range("G9").Value = -0.001110103134557975
$ws.Range("H9").Value = 79.79112358263946
$ws.Range("G10").Value = -0.06162156782579986
$ws.Range("H10").Value = 2.221751575730556
$ws.Range("G11").Value = -0.06757750333692719
$ws.Range("H11").Value = -5.461239003145542
$ws.Range("G12").Value = -0.392819015883466
$ws.Range("H12").Value = 0.4521970569210404
$ws.Range("G13").Value = -0.4101209682551671
$ws.Range("H13").Value = -4.590891776881515
$ws.Range("G14").Value = -0.02547496541445629
$ws.Range("H14").Value = -210.8877981025803
$ws.Range("G15").Value = -0.01814609052928595
$ws.Range("H15").Value = 59.95891242241967
$ws.Range("G16").Value = 0.1391457740419385
$ws.Range("H16").Value = 1.752198061947348
$ws.Range("G17").Value = 0.1395459391505189
$ws.Range("H17").Value = 0.05399031316349422
$ws.Range("G18").Value = 0.1281484326548413
$ws.Range("H18").Value = 8.906590000748876
$ws.Range("G19").Value = 0.1265656535905688
$ws.Range("H19").Value = -1.634581816134293
$ws.Range("G20").Value = 0.08235920787833834
$ws.Range("H20").Value = -7.184821329674042
$ws.Range("G21").Value = 0.08577644863445129
$ws.Range("H21").Value = -1.509428358341044
$ws.Range("G22").Value = -0.1006984782223957
$ws.Range("H22").Value = -7.718297166687712
$ws.Range("G23").Value = -0.1015784942657152
$ws.Range("H23").Value = -0.1315748608146956
$ws.Range("G24").Value = 0.1662150694387406
$ws.Range("H24").Value = 3.180814568098048
$ws.Range("G25").Value = 0.1699360628377588
$ws.Range("H25").Value = -0.3866855973671305
$ws.Range("G26").Value = 0.08715594120014292
$ws.Range("H26").Value = -3.858475491462665
$ws.Range("G27").Value = 0.08792907806349864
$ws.Range("H27").Value = 2.288613097165947
$ws.Range("G28").Value = -0.1349517980561659
$ws.Range("H28").Value = 1.953378219390355
$ws.Range("G29").Value = -0.1420691564941845
$ws.Range("H29").Value = -1.623315330921616
$ws.Range("G30").Value = 0.05240863308751365
$ws.Range("H30").Value = 0.7562365113840086
$ws.Range("G31").Value = 0.0477011399336254
$ws.Range("H31").Value = 8.858488828648802
$ws.Range("G32").Value = 0.107497116889967
$ws.Range("H32").Value = -1.110185029035583
$ws.Range("G33").Value = 0.1125697423380556
$ws.Range("H33").Value = -9.287587918111367
$ws.Range("G34").Value = -0.01299705123635513
$ws.Range("H34").Value = 16.77303172240761
$ws.Range("G35").Value = -0.01802620691242656
$ws.Range("H35").Value = -7.700880058509007
$ws.Range("G36").Value = 0.03130792854008258
$ws.Range("H36").Value = -14.84753832391459
$ws.Range("G37").Value = 0.04008925584117872
$ws.Range("H37").Value = 12.33999539757828
$ws.Range("G38").Value = 0.1013978765757192
$ws.Range("H38").Value = 1.092345810884516
$ws.Range("G39").Value = 0.1006938282631171
$ws.Range("H39").Value = 3.381737975582362
$ws.Range("G40").Value = 0.03059886667688213
$ws.Range("H40").Value = -9.168807690731629
$ws.Range("G41").Value = 0.02625604845921981
$ws.Range("H41").Value = -18.50686695012267
$ws.Range("G42").Value = 0.1209285338811033
$ws.Range("H42").Value = 0.0191086501189757
$ws.Range("G43").Value = 0.1266551918341156
$ws.Range("H43").Value = -0.884494596814097
$ws.Range("G44").Value = 0.04217150620522603
$ws.Range("H44").Value = 6.314795098565834
$ws.Range("G45").Value = 0.03595875262900555
$ws.Range("H45").Value = 15.39164808377636
$ws.Range("G46").Value = 0.05686179819159423
$ws.Range("H46").Value = 0.4311961306809745
$ws.Range("G47").Value = 0.06224557545466997
$ws.Range("H47").Value = 6.094047071282885
$ws.Range("G48").Value = 0.04623503830400241
$ws.Range("H48").Value = -6.13570582273567
$ws.Range("G49").Value = 0.03881861603957311
$ws.Range("H49").Value = -14.81120231732891
$ws.Range("G50").Value = 0.02640027500007253
$ws.Range("H50").Value = -0.3254820620014274
$ws.Range("G51").Value = 0.02837092203851672
$ws.Range("H51").Value = 1.267330259724016
$ws.Range("G52").Value = -0.08625186371179056
$ws.Range("H52").Value = 0.7728362396615259
$ws.Range("G53").Value = -0.07801219005127261
$ws.Range("H53").Value = 2.747589420585881
$ws.Range("G54").Value = 0.04783999368801759
$ws.Range("H54").Value = -4.37108215224801
$ws.Range("G55").Value = 0.05289278510971154
$ws.Range("H55").Value = -6.026949002110161
$ws.Range("G56").Value = 0.04984459331968557
$ws.Range("H56").Value = 0.8348826853728268
$ws.Range("G57").Value = 0.0504877129872548
$ws.Range("H57").Value = 32.91844570833489
$ws.Range("G58").Value = 0.05800623984838288
$ws.Range("H58").Value = 0.6807145494475143
$ws.Range("G59").Value = 0.06340563821663482
$ws.Range("H59").Value = 11.21199420815292
$ws.Range("G60").Value = 0.03795475756637268
$ws.Range("H60").Value = 38.21834753082236
$ws.Range("G61").Value = 0.03481304914890532
$ws.Range("H61").Value = 30.39657251397406
$ws.Range("G62").Value = 0.06323213624758801
$ws.Range("H62").Value = 1.256909169089657
$ws.Range("G63").Value = 0.06453433221337211
$ws.Range("H63").Value = 1.020093345770867
$ws.Range("G64").Value = 0.02250437337694905
$ws.Range("H64").Value = -18.87292249776137
$ws.Range("G65").Value = 0.02886771381773781
$ws.Range("H65").Value = -18.51551486700427
$ws.Range("G66").Value = 0.07937964866103574
$ws.Range("H66").Value = 2.183828558149019
$ws.Range("G67").Value = 0.0822846333244569
$ws.Range("H67").Value = 4.328342984235376
$ws.Range("G68").Value = -0.0230784282516971
$ws.Range("H68").Value = -6.128733329377666
$ws.Range("G69").Value = -0.02500571335017782
$ws.Range("H69").Value = -30.63740402759669
$ws.Range("G70").Value = 0.07555875802645712
$ws.Range("H70").Value = 4.963328987826384
$ws.Range("G71").Value = 0.07020440325943594
$ws.Range("H71").Value = -11.60477731968855
$ws.Range("G72").Value = -0.1498602767390565
$ws.Range("H72").Value = 2.461057753470293
$ws.Range("G73").Value = -0.1514292057344177
$ws.Range("H73").Value = 1.075514051119194
$ws.Range("G74").Value = 0.1531379152630579
$ws.Range("H74").Value = 1.8074298386065
$ws.Range("G75").Value = 0.1496555861218068
$ws.Range("H75").Value = -0.5262562292787016
$ws.Range("G76").Value = -0.001736312633021055
$ws.Range("H76").Value = -67.52309329160883
$ws.Range("G77").Value = -0.005859162608341979
$ws.Range("H77").Value = -165.3741122224224
$ws.Range("G78").Value = 0.09266395776665452
$ws.Range("H78").Value = 3.007736741836759
$ws.Range("G79").Value = 0.09293278869200707
$ws.Range("H79").Value = -4.095480021815605
$ws.Range("G80").Value = -0.2104528402323944
$ws.Range("H80").Value = 2.763218856710008
$ws.Range("G81").Value = -0.2011668992220656
$ws.Range("H81").Value = 5.605916340687964
$ws.Range("G82").Value = 0.1733460959593346
$ws.Range("H82").Value = 3.417556046474623
$ws.Range("G83").Value = 0.1670646425038058
$ws.Range("H83").Value = -5.092509265172435
$ws.Range("G84").Value = 0.09986064149770602
$ws.Range("H84").Value = -5.906975798277764
$ws.Range("G85").Value = 0.1158676633681304
$ws.Range("H85").Value = 10.81578567897644
